$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.360.15'
$ws.Range('E2').Value = '  +0.46%  '
$ws.Range('D3').Value = '2.644.89'
$ws.Range('E3').Value = '  +0.67%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.78'
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.85'
$ws.Range('E6').Value = '  +0.84%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = '2.644.77'
$ws.Range('E9').Value = '  +0.69%  '
$ws.Range('E10').Value = '  +7.43%  '
$ws.Range('E11').Value = '  -0.63%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.27'
$ws.Range('E12').Value = '  +1.19%  '
$ws.Range('E13').Value = '  +1.76%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.23'
$ws.Range('E14').Value = '  +2.02%  '
$ws.Range('E15').Value = '  +2.49%  '
$ws.Range('D16').Value = '3.127.26'
$ws.Range('E16').Value = '  +0.64%  '
$ws.Range('D17').Value = '68.312.07'
$ws.Range('E17').Value = '  +0.62%  '
$ws.Range('D18').Value = '2.638.35'
$ws.Range('E18').Value = '  +0.92%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.44'
$ws.Range('E19').Value = '  +1.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '365.06'
$ws.Range('E20').Value = '  -2.56%  '
$ws.Range('E21').Value = '  +0.79%  '
$ws.Range('E22').Value = '  +3.47%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.90'
$ws.Range('E23').Value = '  +1.62%  '
$ws.Range('E24').Value = '  +0.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '74.49'
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.80'
$ws.Range('E27').Value = '  -1.40%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0000106'
$ws.Range('E28').Value = '  +1.42%  '
$ws.Range('D29').Value = '2.775.41'
$ws.Range('E30').Value = '  -0.16%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '572.58'
$ws.Range('E31').Value = '  -1.00%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.18'
$ws.Range('E32').Value = '  +4.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.42'
$ws.Range('E33').Value = '  +1.30%  '
$ws.Range('E34').Value = '  +0.79%  '
$ws.Range('E35').Value = '  +3.60%  '
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.59'
$ws.Range('E37').Value = '  +5.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '160.85'
$ws.Range('E38').Value = '  +1.28%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.41'
$ws.Range('E39').Value = '  +1.29%  '
$ws.Range('E40').Value = '  +0.33%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.374'
$ws.Range('E41').Value = '  +1.45%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.44'
$ws.Range('E42').Value = '  +1.88%  '
$ws.Range('D43').Value = '0.0₆0338'
$ws.Range('E43').Value = '  +6.39%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.66'
$ws.Range('E44').Value = '  +0.73%  '
$ws.Range('E45').Value = '  +3.60%  '
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '40.47'
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '156.95'
$ws.Range('E48').Value = '  +0.96%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.72'
$ws.Range('E50').Value = '  +1.30%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '21.95'
$ws.Range('E51').Value = '  +0.64%  '
